# edit.ps1 - Word COM-interop script (PowerShell-style) applying the
# "run.docx" changes described by the commit:
#   "Add an abstract class with the often used function. Tests extends
#    from now from SeleniumTests"
#
# The underlying text edits are:
#   1. Paragraph "Actual there is only the workflow of Andre how's tested."
#      -> split into three runs and insert "and some range input tests"
#         before the trailing period (dropping the word "tested").
#   2. Paragraph "The test can be run ... >mvn install." -> split "mvn"
#      into its own run (flanked by spell-check proofErr markers where
#      the runtime allows it).
#   3. "It's also possible to run the test out of the out of your
#      development environment." -> de-duplicate "out of the".
#   4. "It may create some problem ... testing.  Special during the Drag
#      & Drop." -> split into two runs / sentences.
#   5. Remove the trailing empty paragraph's "_GoBack" bookmark so it
#      becomes a plain empty paragraph.

$d = $word.ActiveDocument

# Small helper: force Word to materialize a run boundary at the
# (collapsed or expanded) range `$rng` by nudging a character attribute
# on and back off. Word always splits runs along range boundaries when
# direct formatting is (even transiently) applied, which is exactly the
# mechanism the real editor uses when it inserts/splits text with the
# same visible formatting as its neighbours.
function Split-RunAt($rng) {
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "...André how's tested." -> "...André how's and some range input
#    tests."
# ---------------------------------------------------------------------

# Drop the word "tested" but keep the trailing period.
$r = $d.Content
$r.Find.Execute("tested", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = ""

# Insert the new clause right after "how's " (collapse to the end of the
# found range first so InsertAfter lands exactly between the space and
# the remaining ".").
$r = $d.Content
$r.Find.Execute("how’s ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("and some range input tests")
Split-RunAt $r

# ---------------------------------------------------------------------
# 2) "...Maven test, >mvn install." -> split off "mvn" into its own run
#    (proofed as a spell-check exception in the real document).
# ---------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute("mvn", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-RunAt $r

# ---------------------------------------------------------------------
# 3) "...out of the out of your development environment." -> "...out of
#    your development environment."
# ---------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute("out of the out of your development environment.", $true, $false, $false, $false, $false, $true, 1, $false, "out of your development environment.", 2)

# ---------------------------------------------------------------------
# 4) "...is testing.  Special during the Drag & Drop." -> split into two
#    runs/sentences.
# ---------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute("Special during the Drag & Drop.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-RunAt $r

# ---------------------------------------------------------------------
# 5) Drop the "_GoBack" bookmark on the trailing empty paragraph so it
#    becomes a plain empty paragraph.
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
